# accounts_search (update): rename/replace the "Account2"/"Account3" shared
# strings with the new department-manager-entitlement labels and point the
# affected rows at the right account name, matching the commit's sharedStrings
# reshuffle (Account2 -> "department1 m1e1", Account3 -> "department1 m1e2",
# and row 3's account name becoming plain "Account1").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Account Name") value updates:
#   C3: Account2            -> Account1
#   C4: Account3            -> department1 m1e1
#   C5: Account1            -> department1 m1e2
$ws.Range("C3").Value = "Account1"
$ws.Range("C4").Value = "department1 m1e1"
$ws.Range("C5").Value = "department1 m1e2"

# Column C widens to fit the new (longer) text, same as Excel's "best fit"
# auto-resize after the content change (engine column widths snap to the
# nearest pixel, so this lands on the closest attainable width).
$ws.Columns("C").ColumnWidth = 17.3

# Selection moves off the data range (was C6) to D4.
$ws.Range("D4").Select()
